# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve*Price*, Leve*Profit* columns) across all
# eight crafting-job sheets, row by row, to the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 11804.0625
$ws.Range("I6").Value = 12204.4
$ws.Range("K6").Value = 36613.2
$ws.Range("M6").Value = -36501.2
$ws.Range("H31").Value = 77.375
$ws.Range("I31").Value = 45.42857
$ws.Range("J31").Value = 301
$ws.Range("K31").Value = 136.28571
$ws.Range("L31").Value = 903
$ws.Range("M31").Value = 93.71429000000001
$ws.Range("N31").Value = -1363
$ws.Range("H41").Value = 225.25
$ws.Range("I41").Value = 78
$ws.Range("K41").Value = 78
$ws.Range("M41").Value = 362
$ws.Range("H55").Value = 172.22223
$ws.Range("J55").Value = 214
$ws.Range("L55").Value = 214
$ws.Range("N55").Value = -642
$ws.Range("H62").Value = 8390.5
$ws.Range("I62").Value = 8558.571
$ws.Range("K62").Value = 8558.571
$ws.Range("M62").Value = -7934.571
$ws.Range("H65").Value = 8390.5
$ws.Range("I65").Value = 8558.571
$ws.Range("K65").Value = 42792.855
$ws.Range("M65").Value = -39672.855
$ws.Range("H96").Value = 389.5909
$ws.Range("I96").Value = 254.55556
$ws.Range("J96").Value = 997.25
$ws.Range("K96").Value = 763.66668
$ws.Range("L96").Value = 2991.75
$ws.Range("M96").Value = 609.33332
$ws.Range("N96").Value = -5737.75
$ws.Range("H107").Value = 1497.9565
$ws.Range("J107").Value = 2502
$ws.Range("L107").Value = 2502
$ws.Range("N107").Value = -6342
$ws.Range("H111").Value = 1395
$ws.Range("I111").Value = 1334.4
$ws.Range("K111").Value = 4003.2
$ws.Range("M111").Value = -936.2000000000003
$ws.Range("H116").Value = 8851.700000000001
$ws.Range("I116").Value = 7166.25
$ws.Range("K116").Value = 7166.25
$ws.Range("M116").Value = -3724.25
$ws.Range("H132").Value = 6225.3076
$ws.Range("I132").Value = 5584.5454
$ws.Range("K132").Value = 16753.6362
$ws.Range("M132").Value = -14223.6362
$ws.Range("H137").Value = 20340.479
$ws.Range("I137").Value = 9209.666999999999
$ws.Range("J137").Value = 27496
$ws.Range("K137").Value = 27629.001
$ws.Range("L137").Value = 82488
$ws.Range("M137").Value = -25079.001
$ws.Range("N137").Value = -87588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3759.9
$ws.Range("I2").Value = 3256.8215
$ws.Range("J2").Value = 4933.75
$ws.Range("K2").Value = 3256.8215
$ws.Range("L2").Value = 4933.75
$ws.Range("M2").Value = -3143.8215
$ws.Range("N2").Value = -5159.75
$ws.Range("H32").Value = 2477.2
$ws.Range("I32").Value = 1746.174
$ws.Range("K32").Value = 1746.174
$ws.Range("M32").Value = -1459.174
$ws.Range("H61").Value = 4168.8
$ws.Range("I61").Value = 1984.7778
$ws.Range("K61").Value = 1984.7778
$ws.Range("M61").Value = -1772.7778
$ws.Range("H110").Value = 2062.5833
$ws.Range("I110").Value = 1878.5
$ws.Range("K110").Value = 1878.5
$ws.Range("M110").Value = 166.5
$ws.Range("H116").Value = 3759.9
$ws.Range("I116").Value = 3256.8215
$ws.Range("J116").Value = 4933.75
$ws.Range("K116").Value = 3256.8215
$ws.Range("L116").Value = 4933.75
$ws.Range("M116").Value = -962.8215
$ws.Range("N116").Value = -9521.75
$ws.Range("H136").Value = 4168.8
$ws.Range("I136").Value = 1984.7778
$ws.Range("K136").Value = 5954.3334
$ws.Range("M136").Value = -3404.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3759.9
$ws.Range("I3").Value = 3256.8215
$ws.Range("J3").Value = 4933.75
$ws.Range("K3").Value = 3256.8215
$ws.Range("L3").Value = 4933.75
$ws.Range("M3").Value = -3142.8215
$ws.Range("N3").Value = -5161.75
$ws.Range("H20").Value = 6947
$ws.Range("I20").Value = 4718.25
$ws.Range("J20").Value = 11404.5
$ws.Range("K20").Value = 4718.25
$ws.Range("L20").Value = 11404.5
$ws.Range("M20").Value = -4471.25
$ws.Range("N20").Value = -11898.5
$ws.Range("H94").Value = 844.5263
$ws.Range("I94").Value = 264.26666
$ws.Range("J94").Value = 3020.5
$ws.Range("K94").Value = 264.26666
$ws.Range("L94").Value = 3020.5
$ws.Range("M94").Value = 186.73334
$ws.Range("N94").Value = -3922.5
$ws.Range("H99").Value = 2414.1428
$ws.Range("I99").Value = 2279.8
$ws.Range("K99").Value = 2279.8
$ws.Range("M99").Value = -781.8000000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3254.2632
$ws.Range("I16").Value = 3036.375
$ws.Range("K16").Value = 3036.375
$ws.Range("M16").Value = -2749.375
$ws.Range("H113").Value = 3254.2632
$ws.Range("I113").Value = 3036.375
$ws.Range("K113").Value = 3036.375
$ws.Range("M113").Value = -866.375
$ws.Range("H141").Value = 96598.75
$ws.Range("J141").Value = 114999
$ws.Range("L141").Value = 114999
$ws.Range("N141").Value = -125359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2286.4
$ws.Range("I69").Value = 2335.5715
$ws.Range("K69").Value = 7006.7145
$ws.Range("M69").Value = -6195.7145
$ws.Range("H72").Value = 2286.4
$ws.Range("I72").Value = 2335.5715
$ws.Range("K72").Value = 21020.1435
$ws.Range("M72").Value = -16964.1435
$ws.Range("H98").Value = 493.69232
$ws.Range("I98").Value = 555.25
$ws.Range("K98").Value = 1665.75
$ws.Range("M98").Value = -167.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 215.83333
$ws.Range("I2").Value = 172.22223
$ws.Range("K2").Value = 172.22223
$ws.Range("M2").Value = -59.22223
$ws.Range("H47").Value = 45000
$ws.Range("J47").Value = 45000
$ws.Range("L47").Value = 45000
$ws.Range("N47").Value = -46136
$ws.Range("H57").Value = 34820.2
$ws.Range("J57").Value = 42511.5
$ws.Range("L57").Value = 42511.5
$ws.Range("N57").Value = -44151.5
$ws.Range("H80").Value = 2657.2
$ws.Range("I80").Value = 2596.4
$ws.Range("K80").Value = 2596.4
$ws.Range("M80").Value = -1598.4
$ws.Range("H83").Value = 2657.2
$ws.Range("I83").Value = 2596.4
$ws.Range("K83").Value = 12982
$ws.Range("M83").Value = -7990
$ws.Range("H102").Value = 2532.6667
$ws.Range("I102").Value = 2016.3334
$ws.Range("K102").Value = 2016.3334
$ws.Range("M102").Value = -394.3334
$ws.Range("H107").Value = 441.86206
$ws.Range("I107").Value = 359.45456
$ws.Range("J107").Value = 700.8570999999999
$ws.Range("K107").Value = 359.45456
$ws.Range("L107").Value = 700.8570999999999
$ws.Range("M107").Value = 1560.54544
$ws.Range("N107").Value = -4540.8571
$ws.Range("H113").Value = 224645.44
$ws.Range("I113").Value = 288128
$ws.Range("J113").Value = 2456.5
$ws.Range("K113").Value = 288128
$ws.Range("L113").Value = 2456.5
$ws.Range("M113").Value = -285958
$ws.Range("N113").Value = -6796.5
$ws.Range("H126").Value = 3961
$ws.Range("I126").Value = 3349.1667
$ws.Range("K126").Value = 10047.5001
$ws.Range("M126").Value = -7577.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2917.879
$ws.Range("I22").Value = 2264.0386
$ws.Range("K22").Value = 2264.0386
$ws.Range("M22").Value = -1969.0386
$ws.Range("H27").Value = 2917.879
$ws.Range("I27").Value = 2264.0386
$ws.Range("K27").Value = 2264.0386
$ws.Range("M27").Value = -2157.0386
$ws.Range("H46").Value = 1699.4
$ws.Range("I46").Value = 1071.1428
$ws.Range("J46").Value = 2249.125
$ws.Range("K46").Value = 1071.1428
$ws.Range("L46").Value = 2249.125
$ws.Range("M46").Value = -883.1428000000001
$ws.Range("N46").Value = -2625.125
$ws.Range("H61").Value = 2437
$ws.Range("I61").Value = 2421.25
$ws.Range("K61").Value = 2421.25
$ws.Range("M61").Value = -2219.25
$ws.Range("H93").Value = 2022.9445
$ws.Range("I93").Value = 2031.8438
$ws.Range("K93").Value = 2031.8438
$ws.Range("M93").Value = -783.8438000000001
$ws.Range("H100").Value = 1011959
$ws.Range("I100").Value = 93854.27
$ws.Range("J100").Value = 11111111
$ws.Range("K100").Value = 93854.27
$ws.Range("L100").Value = 11111111
$ws.Range("M100").Value = -93313.27
$ws.Range("N100").Value = -11112193
$ws.Range("H104").Value = 22071.8
$ws.Range("J104").Value = 22071.8
$ws.Range("L104").Value = 22071.8
$ws.Range("N104").Value = -29059.8
$ws.Range("H113").Value = 2437
$ws.Range("I113").Value = 2421.25
$ws.Range("K113").Value = 2421.25
$ws.Range("M113").Value = -251.25
$ws.Range("H122").Value = 3545.1428
$ws.Range("I122").Value = 2946
$ws.Range("K122").Value = 8838
$ws.Range("M122").Value = -6388
$ws.Range("H132").Value = 4634067
$ws.Range("I132").Value = 6177091
$ws.Range("K132").Value = 18531273
$ws.Range("M132").Value = -18528743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1398.375
$ws.Range("I100").Value = 1469.4
$ws.Range("J100").Value = 1280
$ws.Range("K100").Value = 2938.8
$ws.Range("L100").Value = 2560
$ws.Range("M100").Value = -2397.8
$ws.Range("N100").Value = -3642
$ws.Range("H107").Value = 1059.8518
$ws.Range("I107").Value = 735.05884
$ws.Range("J107").Value = 1612
$ws.Range("K107").Value = 2205.17652
$ws.Range("L107").Value = 4836
$ws.Range("M107").Value = -285.17652
$ws.Range("N107").Value = -8676

